# Append newly-studied vocabulary entries to the ENGLISH sheet and newly
# collected quotes/notes to the NOTES sheet, tracking everything added
# during this session.

$wb = $excel.ActiveWorkbook
$english = $wb.Worksheets.Item("ENGLISH")
$notes = $wb.Worksheets.Item("NOTES")

# --- New vocabulary words appended to ENGLISH (columns: Word, Definition,
#     Synonyms, Antonyms, Correct answer count, Created at, Tags) ---
$words = @(
    @("congested", "", "crowded", "", 0, "2021-11-20 22:32:31.286541", ""),
    @("swerve", "change or cause to change direction abruptly", "veer", "", 0, "2021-11-20 22:33:54.608594", ""),
    @("magnum opus", "a work of art, music, or literature that is regarded as the most important or best work that an artist, composer, or writer has produced", "", "", 0, "2021-11-20 22:36:19.119537", ""),
    @("anticipate", "", "expect;predict", "", 0, "2021-11-20 22:36:48.436207", ""),
    @("windfall", "a large amount of money that is won or received unexpectedly", "", "", 0, "2021-11-20 22:52:03.527958", ""),
    @("cushion", "", "pillow;protection", "", 0, "2021-11-20 22:53:43.81857", ""),
    @("diligently", "in a way that shows care in one's work or duties", "", "", 0, "2021-11-20 22:55:00.4924", ""),
    @("resilience", "", "flexibility", "", 0, "2021-11-20 22:56:23.511821", ""),
    @("parable", "", "allegory", "", 0, "2021-11-20 22:57:34.943717", ""),
    @("arbitrary", "", "random", "", 0, "2021-11-20 22:59:11.038374", "")
)

$startRow = $english.Cells.Item($english.Rows.Count, 1).End(-4162).Row + 1
$r = $startRow
foreach ($word in $words) {
    $english.Cells.Item($r, 1).Value = $word[0]
    $english.Cells.Item($r, 2).Value = $word[1]
    $english.Cells.Item($r, 3).Value = $word[2]
    $english.Cells.Item($r, 4).Value = $word[3]
    $english.Cells.Item($r, 5).Value = $word[4]
    $english.Cells.Item($r, 6).Value = $word[5]
    $english.Cells.Item($r, 7).Value = $word[6]
    $r = $r + 1
}

# --- New quotes/notes appended to NOTES (columns: Note, Tags) ---
$notesToAdd = @(
    "Give me 6 hours to chop down a tree and I will spend the first 4 sharpening the axe",
    "The only thing we can expect (with any great certainty) is the unexpected",
    "To attain knowledge add things every day. To attain wisdom subtract things every day"
)

$noteStartRow = $notes.Cells.Item($notes.Rows.Count, 1).End(-4162).Row + 1
$r = $noteStartRow
foreach ($note in $notesToAdd) {
    $notes.Cells.Item($r, 1).Value = $note
    $notes.Cells.Item($r, 2).Value = ""
    $r = $r + 1
}
